$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New scheme names that now occupy rows 10-13 (Gaussian-Quadrature + 3 spiral
# schemes), pushing the previously-existing schemes (NoRotation-tilt60deg,
# Rotation-NoTilt, Rotation-60detTilt, HexGrid-*) down by 4 rows. Three of
# those displaced schemes spill into brand-new rows 17-19.
$schemeNames = @(
    "ND Single",
    "RD Single",
    "TD Single",
    "Morris",
    "Ring Perpendicular to ND",
    "Ring Perpendicular to RD",
    "Ring Perpendicular to TD",
    "Gaussian-Quadrature",
    "Spiral-90deg-10rot-5space",
    "Spiral-90deg-15rot-5space",
    "Spiral-90deg-10rot-3space",
    "NoRotation-tilt60deg",
    "Rotation-NoTilt",
    "Rotation-60detTilt",
    "HexGrid-90degTilt5degRes",
    "HexGrid-90degTilt22p5degRes",
    "HexGrid-60degTilt5degRes"
)

for ($i = 0; $i -lt $schemeNames.Length; $i++) {
    $row = $i + 3
    $ws.Cells.Item($row, 1).Value = $i + 1
    $ws.Cells.Item($row, 2).Value = $schemeNames[$i]
    for ($col = 3; $col -le 13; $col++) {
        $ws.Cells.Item($row, $col).Value = 1
    }
}

# The 3 brand-new rows (17-19) need the same "index column" formatting
# (bold / centered / bordered) that column A already carries on every
# existing data row. Copy that formatting down from the last pre-existing
# row (A16) onto the newly created rows.
$ws.Range("A16").Copy() | Out-Null
$ws.Range("A17:A19").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = 0
